$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to match the refreshed crypto data.
# Set NumberFormat to text ("@") per cell first so that numeric-looking strings
# (e.g. "1.001", "0.00001134") are preserved exactly as text, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.387.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.569.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.89"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3678"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.31"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3378"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.164"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07619"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.27"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.055"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.913"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.577.15"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001134"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.64"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06751"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.219"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5310"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.49"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.71%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.396.59"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.365"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.911"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.98"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "146.03"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.983"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.41"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.741.43"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.056"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +8.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.249"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.009"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.23"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -9.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08459"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02546"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.324"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.69"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -8.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6346"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.19"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.45%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5973"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.752"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.103"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.260"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.71"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.66%  "

# Rows 40 and 41 swap places: "InternetComputer(DFINITY)" moves to row 41
# and "Hedera" moves to row 40, each carrying freshly updated price/volume data.
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06527"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.522"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.76%  "
